$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$passwords = @{
    2  = "crEVBCu"
    3  = "DzLQJXF"
    4  = "xwC6c09"
    5  = "D5xeCHW"
    6  = "HXk81eY"
    7  = "BA3nyT4"
    8  = "qj7uAMD"
    9  = "Le17bMO"
    10 = "NdvO6NV"
    11 = "GwPMOZd"
    12 = "uXZduwT"
    13 = "a8FQgNH"
    14 = "TFVst32"
    15 = "AF9BVvx"
    16 = "IPialY3"
    17 = "qhz0yx8"
    18 = "B2kSZUD"
    19 = "CUPusvA"
    20 = "28Z0y8V"
    21 = "vB9hiWd"
    22 = "7C3b7XL"
    23 = "zCRPpAZ"
    24 = "xKkY4x6"
    25 = "PwV5ONe"
    26 = "nrxXtoX"
    27 = "syS0VCN"
    28 = "WeFGbxw"
    29 = "4hivFRY"
    30 = "3dxpPOy"
    31 = "JaXJJBK"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 4).Value = $passwords[$row]
}
